# Auto-generated edit script: applies betexplorer re-scrape diff to germany_bundesliga_2023-2024 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose match data (columns F:V) shifted position after re-scrape;
#     identity columns A:E (Indice/pais/torneio/temporada/data_partida) stay put. ---
# Row 51
$ws.Cells.Item(51, 6).Value = "Bochum"
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(51, 8).Value = "B. Monchengladbach"
$ws.Cells.Item(51, 9).Value = 3
$ws.Cells.Item(51, 10).Value = 2.17
$ws.Cells.Item(51, 11).Value = "17/09/2023 09:01"
$ws.Cells.Item(51, 12).Value = 2.27
$ws.Cells.Item(51, 13).Value = "30/09/2023 15:27"
$ws.Cells.Item(51, 14).Value = 3.91
$ws.Cells.Item(51, 15).Value = "17/09/2023 09:01"
$ws.Cells.Item(51, 16).Value = 3.77
$ws.Cells.Item(51, 17).Value = "30/09/2023 15:27"
$ws.Cells.Item(51, 18).Value = 3.27
$ws.Cells.Item(51, 19).Value = "17/09/2023 09:01"
$ws.Cells.Item(51, 20).Value = 3.15
$ws.Cells.Item(51, 21).Value = "30/09/2023 15:26"
$ws.Cells.Item(51, 22).Value = "https://www.betexplorer.com/football/germany/bundesliga/bochum-b-monchengladbach/fwUFeisg/"

# Row 52
$ws.Cells.Item(52, 6).Value = "FC Koln"
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = "Stuttgart"
$ws.Cells.Item(52, 9).Value = 2
$ws.Cells.Item(52, 10).Value = 2.49
$ws.Cells.Item(52, 11).Value = "17/09/2023 09:01"
$ws.Cells.Item(52, 12).Value = 2.86
$ws.Cells.Item(52, 13).Value = "30/09/2023 15:23"
$ws.Cells.Item(52, 14).Value = 3.59
$ws.Cells.Item(52, 15).Value = "17/09/2023 09:01"
$ws.Cells.Item(52, 16).Value = 3.58
$ws.Cells.Item(52, 17).Value = "30/09/2023 15:20"
$ws.Cells.Item(52, 18).Value = 2.92
$ws.Cells.Item(52, 19).Value = "17/09/2023 09:01"
$ws.Cells.Item(52, 20).Value = 2.55
$ws.Cells.Item(52, 21).Value = "30/09/2023 15:20"
$ws.Cells.Item(52, 22).Value = "https://www.betexplorer.com/football/germany/bundesliga/1-fc-koln-vfb-stuttgart/nDQ7cDCt/"

# Row 57
$ws.Cells.Item(57, 6).Value = "Augsburg"
$ws.Cells.Item(57, 7).Value = 1
$ws.Cells.Item(57, 8).Value = "Darmstadt"
$ws.Cells.Item(57, 9).Value = 2
$ws.Cells.Item(57, 10).Value = 1.8
$ws.Cells.Item(57, 11).Value = "28/09/2023 14:20"
$ws.Cells.Item(57, 12).Value = 1.85
$ws.Cells.Item(57, 13).Value = "07/10/2023 15:01"
$ws.Cells.Item(57, 14).Value = 3.89
$ws.Cells.Item(57, 15).Value = "28/09/2023 14:20"
$ws.Cells.Item(57, 16).Value = 3.75
$ws.Cells.Item(57, 17).Value = "07/10/2023 15:29"
$ws.Cells.Item(57, 18).Value = 4.06
$ws.Cells.Item(57, 19).Value = "28/09/2023 14:20"
$ws.Cells.Item(57, 20).Value = 4.66
$ws.Cells.Item(57, 21).Value = "07/10/2023 15:28"
$ws.Cells.Item(57, 22).Value = "https://www.betexplorer.com/football/germany/bundesliga/augsburg-darmstadt/fa5j8UQb/"

# Row 58
$ws.Cells.Item(58, 6).Value = "Dortmund"
$ws.Cells.Item(58, 7).Value = 4
$ws.Cells.Item(58, 8).Value = "Union Berlin"
$ws.Cells.Item(58, 9).Value = 2
$ws.Cells.Item(58, 10).Value = 1.52
$ws.Cells.Item(58, 11).Value = "23/09/2023 19:02"
$ws.Cells.Item(58, 12).Value = 1.71
$ws.Cells.Item(58, 13).Value = "07/10/2023 15:25"
$ws.Cells.Item(58, 14).Value = 4.52
$ws.Cells.Item(58, 15).Value = "23/09/2023 19:02"
$ws.Cells.Item(58, 16).Value = 4.17
$ws.Cells.Item(58, 17).Value = "07/10/2023 15:27"
$ws.Cells.Item(58, 18).Value = 6.35
$ws.Cells.Item(58, 19).Value = "23/09/2023 19:02"
$ws.Cells.Item(58, 20).Value = 4.97
$ws.Cells.Item(58, 21).Value = "07/10/2023 15:29"
$ws.Cells.Item(58, 22).Value = "https://www.betexplorer.com/football/germany/bundesliga/dortmund-union-berlin/OILWi9tI/"

# Row 59
$ws.Cells.Item(59, 6).Value = "RB Leipzig"
$ws.Cells.Item(59, 7).Value = 0
$ws.Cells.Item(59, 8).Value = "Bochum"
$ws.Cells.Item(59, 9).Value = 0
$ws.Cells.Item(59, 10).Value = 1.31
$ws.Cells.Item(59, 11).Value = "23/09/2023 21:02"
$ws.Cells.Item(59, 12).Value = 1.26
$ws.Cells.Item(59, 13).Value = "07/10/2023 15:20"
$ws.Cells.Item(59, 14).Value = 5.99
$ws.Cells.Item(59, 15).Value = "23/09/2023 21:02"
$ws.Cells.Item(59, 16).Value = 6.66
$ws.Cells.Item(59, 17).Value = "07/10/2023 15:29"
$ws.Cells.Item(59, 18).Value = 9.130000000000001
$ws.Cells.Item(59, 19).Value = "23/09/2023 21:02"
$ws.Cells.Item(59, 20).Value = 11.62
$ws.Cells.Item(59, 21).Value = "07/10/2023 15:29"
$ws.Cells.Item(59, 22).Value = "https://www.betexplorer.com/football/germany/bundesliga/rb-leipzig-bochum/4UxUjTeO/"

# Row 66
$ws.Cells.Item(66, 6).Value = "Freiburg"
$ws.Cells.Item(66, 7).Value = 2
$ws.Cells.Item(66, 8).Value = "Bochum"
$ws.Cells.Item(66, 9).Value = 1
$ws.Cells.Item(66, 10).Value = 1.56
$ws.Cells.Item(66, 11).Value = "01/10/2023 23:02"
$ws.Cells.Item(66, 12).Value = 1.81
$ws.Cells.Item(66, 13).Value = "21/10/2023 15:24"
$ws.Cells.Item(66, 14).Value = 4.28
$ws.Cells.Item(66, 15).Value = "01/10/2023 23:02"
$ws.Cells.Item(66, 16).Value = 3.91
$ws.Cells.Item(66, 17).Value = "21/10/2023 15:27"
$ws.Cells.Item(66, 18).Value = 5.28
$ws.Cells.Item(66, 19).Value = "01/10/2023 23:02"
$ws.Cells.Item(66, 20).Value = 4.61
$ws.Cells.Item(66, 21).Value = "21/10/2023 15:26"
$ws.Cells.Item(66, 22).Value = "https://www.betexplorer.com/football/germany/bundesliga/freiburg-bochum/rZG64nQN/"

# Row 67
$ws.Cells.Item(67, 6).Value = "Hoffenheim"
$ws.Cells.Item(67, 7).Value = 1
$ws.Cells.Item(67, 8).Value = "Eintracht Frankfurt"
$ws.Cells.Item(67, 9).Value = 3
$ws.Cells.Item(67, 10).Value = 2.39
$ws.Cells.Item(67, 11).Value = "01/10/2023 23:01"
$ws.Cells.Item(67, 12).Value = 2.06
$ws.Cells.Item(67, 13).Value = "21/10/2023 15:28"
$ws.Cells.Item(67, 14).Value = 3.56
$ws.Cells.Item(67, 15).Value = "01/10/2023 23:01"
$ws.Cells.Item(67, 16).Value = 3.65
$ws.Cells.Item(67, 17).Value = "21/10/2023 15:26"
$ws.Cells.Item(67, 18).Value = 2.98
$ws.Cells.Item(67, 19).Value = "01/10/2023 23:01"
$ws.Cells.Item(67, 20).Value = 3.79
$ws.Cells.Item(67, 21).Value = "21/10/2023 15:28"
$ws.Cells.Item(67, 22).Value = "https://www.betexplorer.com/football/germany/bundesliga/hoffenheim-eintracht-frankfurt/fDpPMTuh/"

# Row 68
$ws.Cells.Item(68, 6).Value = "Union Berlin"
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = "Stuttgart"
$ws.Cells.Item(68, 9).Value = 3
$ws.Cells.Item(68, 10).Value = 2.3
$ws.Cells.Item(68, 11).Value = "01/10/2023 23:01"
$ws.Cells.Item(68, 12).Value = 2.61
$ws.Cells.Item(68, 13).Value = "21/10/2023 15:29"
$ws.Cells.Item(68, 14).Value = 3.46
$ws.Cells.Item(68, 15).Value = "01/10/2023 23:01"
$ws.Cells.Item(68, 16).Value = 3.37
$ws.Cells.Item(68, 17).Value = "21/10/2023 15:28"
$ws.Cells.Item(68, 18).Value = 3.34
$ws.Cells.Item(68, 19).Value = "01/10/2023 23:01"
$ws.Cells.Item(68, 20).Value = 2.92
$ws.Cells.Item(68, 21).Value = "21/10/2023 15:29"
$ws.Cells.Item(68, 22).Value = "https://www.betexplorer.com/football/germany/bundesliga/union-berlin-vfb-stuttgart/plb35SAH/"

# Row 69
$ws.Cells.Item(69, 6).Value = "Wolfsburg"
$ws.Cells.Item(69, 7).Value = 1
$ws.Cells.Item(69, 8).Value = "Bayer Leverkusen"
$ws.Cells.Item(69, 9).Value = 2
$ws.Cells.Item(69, 10).Value = 2.83
$ws.Cells.Item(69, 11).Value = "01/10/2023 23:01"
$ws.Cells.Item(69, 12).Value = 4.47
$ws.Cells.Item(69, 13).Value = "21/10/2023 15:00"
$ws.Cells.Item(69, 14).Value = 3.69
$ws.Cells.Item(69, 15).Value = "01/10/2023 23:01"
$ws.Cells.Item(69, 16).Value = 4.11
$ws.Cells.Item(69, 17).Value = "21/10/2023 14:58"
$ws.Cells.Item(69, 18).Value = 2.52
$ws.Cells.Item(69, 19).Value = "01/10/2023 23:01"
$ws.Cells.Item(69, 20).Value = 1.79
$ws.Cells.Item(69, 21).Value = "21/10/2023 15:00"
$ws.Cells.Item(69, 22).Value = "https://www.betexplorer.com/football/germany/bundesliga/wolfsburg-bayer-leverkusen/xdGA36uU/"

# Row 70
$ws.Cells.Item(70, 6).Value = "Darmstadt"
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = "RB Leipzig"
$ws.Cells.Item(70, 9).Value = 3
$ws.Cells.Item(70, 10).Value = 6.5
$ws.Cells.Item(70, 11).Value = "02/10/2023 08:32"
$ws.Cells.Item(70, 12).Value = 6.59
$ws.Cells.Item(70, 13).Value = "21/10/2023 15:29"
$ws.Cells.Item(70, 14).Value = 5.15
$ws.Cells.Item(70, 15).Value = "02/10/2023 08:32"
$ws.Cells.Item(70, 16).Value = 5.1
$ws.Cells.Item(70, 17).Value = "21/10/2023 15:29"
$ws.Cells.Item(70, 18).Value = 1.4
$ws.Cells.Item(70, 19).Value = "02/10/2023 08:32"
$ws.Cells.Item(70, 20).Value = 1.47
$ws.Cells.Item(70, 21).Value = "21/10/2023 15:28"
$ws.Cells.Item(70, 22).Value = "https://www.betexplorer.com/football/germany/bundesliga/darmstadt-rb-leipzig/2wlXK7A4/"

# Row 77
$ws.Cells.Item(77, 6).Value = "Augsburg"
$ws.Cells.Item(77, 7).Value = 3
$ws.Cells.Item(77, 8).Value = "Wolfsburg"
$ws.Cells.Item(77, 9).Value = 2
$ws.Cells.Item(77, 10).Value = 2.72
$ws.Cells.Item(77, 11).Value = "10/10/2023 14:02"
$ws.Cells.Item(77, 12).Value = 2.54
$ws.Cells.Item(77, 13).Value = "28/10/2023 14:56"
$ws.Cells.Item(77, 14).Value = 3.62
$ws.Cells.Item(77, 15).Value = "10/10/2023 14:02"
$ws.Cells.Item(77, 16).Value = 3.62
$ws.Cells.Item(77, 17).Value = "28/10/2023 14:56"
$ws.Cells.Item(77, 18).Value = 2.44
$ws.Cells.Item(77, 19).Value = "10/10/2023 14:02"
$ws.Cells.Item(77, 20).Value = 2.8
$ws.Cells.Item(77, 21).Value = "28/10/2023 14:56"
$ws.Cells.Item(77, 22).Value = "https://www.betexplorer.com/football/germany/bundesliga/augsburg-wolfsburg/d4u8MKo9/"

# Row 79
$ws.Cells.Item(79, 6).Value = "Werder Bremen"
$ws.Cells.Item(79, 7).Value = 2
$ws.Cells.Item(79, 8).Value = "Union Berlin"
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 2.73
$ws.Cells.Item(79, 11).Value = "10/10/2023 14:02"
$ws.Cells.Item(79, 12).Value = 2.86
$ws.Cells.Item(79, 13).Value = "28/10/2023 15:01"
$ws.Cells.Item(79, 14).Value = 3.33
$ws.Cells.Item(79, 15).Value = "10/10/2023 14:02"
$ws.Cells.Item(79, 16).Value = 3.53
$ws.Cells.Item(79, 17).Value = "28/10/2023 14:50"
$ws.Cells.Item(79, 18).Value = 2.73
$ws.Cells.Item(79, 19).Value = "10/10/2023 14:02"
$ws.Cells.Item(79, 20).Value = 2.53
$ws.Cells.Item(79, 21).Value = "28/10/2023 15:01"
$ws.Cells.Item(79, 22).Value = "https://www.betexplorer.com/football/germany/bundesliga/werder-bremen-union-berlin/4AhMbHg2/"

# --- New row 101: match added by the 24-11-2023 20:45 re-scrape run ---
# Clone formatting (bold/border/center style for A, date-time style for E) from row 100
# so the new row reuses the existing cell styles instead of creating new ones.
$ws.Cells.Item(100, 1).Copy($ws.Cells.Item(101, 1))
$ws.Cells.Item(100, 5).Copy($ws.Cells.Item(101, 5))

$ws.Cells.Item(101, 1).Value = 100
$ws.Cells.Item(101, 2).Value = "germany"
$ws.Cells.Item(101, 3).Value = "bundesliga"
$ws.Cells.Item(101, 4).Value = "2023-2024"
$ws.Cells.Item(101, 5).Value = 45254.85416666666
$ws.Cells.Item(101, 6).Value = "FC Koln"
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = "Bayern Munich"
$ws.Cells.Item(101, 9).Value = 1
$ws.Cells.Item(101, 10).Value = 6.77
$ws.Cells.Item(101, 11).Value = "05/11/2023 11:03"
$ws.Cells.Item(101, 12).Value = 10.08
$ws.Cells.Item(101, 13).Value = "24/11/2023 20:29"
$ws.Cells.Item(101, 14).Value = 5.6
$ws.Cells.Item(101, 15).Value = "05/11/2023 11:03"
$ws.Cells.Item(101, 16).Value = 6.98
$ws.Cells.Item(101, 17).Value = "24/11/2023 20:29"
$ws.Cells.Item(101, 18).Value = 1.4
$ws.Cells.Item(101, 19).Value = "05/11/2023 11:03"
$ws.Cells.Item(101, 20).Value = 1.26
$ws.Cells.Item(101, 21).Value = "24/11/2023 19:30"
$ws.Cells.Item(101, 22).Value = "https://www.betexplorer.com/football/germany/bundesliga/1-fc-koln-bayern-munich/EwL6NLm0/"